$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F11").Value = 81
$ws1.Range("F13").Value = 1491
$ws1.Range("F15").Value = 39
$ws1.Range("F16").Value = 321
$ws1.Range("G17").Value = 30
$ws1.Range("F23").Value = 5
$ws1.Range("F26").Value = 1551
$ws1.Range("F27").Value = 22
$ws1.Range("F29").Value = 407
$ws1.Range("F32").Value = 395

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 81
$ws4.Range("F14").Value = 1491
$ws4.Range("F16").Value = 39
$ws4.Range("F17").Value = 321
$ws4.Range("G18").Value = 30
$ws4.Range("F24").Value = 5
$ws4.Range("F27").Value = 1551
$ws4.Range("F28").Value = 22
$ws4.Range("F30").Value = 407
$ws4.Range("F33").Value = 395
